$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP value in C2 (mirrors a user typing a new value into the cell)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "192.168.1.113"

# Excel auto-fits the column to the new (longer) content; B keeps its old
# width, C grows to fit "192.168.1.113" (resulting stored width is 15;
# the host adds ~5/7 of a character to whatever ColumnWidth is assigned,
# so back that off here to land exactly on 15)
$ws.Columns("C").ColumnWidth = 14.285714285714286

# The edited cell becomes the active selection
$ws.Range("C2").Select()
